# Add a "Pop" column (column D) to the "income" worksheet.
# Mirrors the commit "Added Pop column to income-data": a new header
# cell "Pop" plus ten population figures, one per existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("income")

# New values for column D: header, then one number per data row (2-11).
$header = "Pop"
$popValues = @(70, 33, 23, 54, 11, 27, 81, 24, 10, 7)

# Cell each new D-cell should copy its formatting from, so the copied
# style reuses an existing style entry instead of creating a new one
# (mirrors the mixed formatting already present on the sheet).
$styleSource = @("C1", "A2", "B3", "A4", "B5", "B6", "A7", "A8", "A9", "A10", "A11")

# Header
$ws.Range($styleSource[0]).Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = $header

# Data rows
for ($i = 0; $i -lt $popValues.Length; $i++) {
    $row = $i + 2
    $target = $ws.Range("D" + $row)
    $ws.Range($styleSource[$i + 1]).Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null
    $target.Value = $popValues[$i]
}

$excel.CutCopyMode = 0
